$d = $word.ActiveDocument

# The document ends with a bulleted ("Prrafodelista") paragraph that carries
# the hidden "_GoBack" bookmark right after its text (before the paragraph
# mark). We need to:
#   1. Add a new bullet paragraph after it with the new comment text.
#   2. Move the "_GoBack" bookmark from the old last paragraph onto the new
#      last paragraph (Word always keeps "_GoBack" at the most recent edit).

$oldLast = $d.Paragraphs.Last

# Append a new paragraph right after the current last one; Word mirrors the
# paragraph's style/numbering (Prrafodelista / numId 2) onto the new one.
$oldLast.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last

# Type the new text plus a one-character sentinel ("X"). Writing the
# sentinel lets us park the bookmark immediately to its *left* (i.e. right
# after the real text) without the bookmark's collapsed range coinciding
# with the paragraph-mark position - doing that directly trips an edge case
# in Bookmarks.Add at end-of-paragraph offsets.
$newLast.Range.Text = "Agregar funcionalidad para agregar comentarios" + "X"
$newLast = $d.Paragraphs.Last

$sentinelStart = $newLast.Range.End - 2
$bookmarkPos = $d.Range($sentinelStart, $sentinelStart)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkPos)

# Remove the sentinel character now that the bookmark is anchored right
# before it; the collapsed bookmark stays put, ending up exactly at the new
# paragraph's end (right after the real text, before the paragraph mark).
$newLast = $d.Paragraphs.Last
$sentinelRange = $d.Range($newLast.Range.End - 2, $newLast.Range.End - 1)
$sentinelRange.Delete()
